$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 32 (محمد رضا مددی هیر): remove the placeholder numeric code entirely.
$ws.Range("C32").ClearContents()

# Append the new row 37 (خانم شیبانی / معاون / she) first so the three new
# shared strings it introduces land before the ones touched below - this
# mirrors the order the strings were appended to sharedStrings.xml.
$ws.Range("A37").Value = "خانم شیبانی"
$ws.Range("B37").Value = "معاون"
$ws.Range("C37").Value = "she"

# Replace the remaining numeric placeholders (123) in column C with their
# real short text codes, in the same order the new shared strings appear.
$ws.Range("C36").Value = "taj"
$ws.Range("C34").Value = "poy"
$ws.Range("C33").Value = "arm"

# Match the author's final selection state.
$ws.Range("C32").Select()
